$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$co = $ws.ChartObjects(1)
$co.Width = 618.125
$co.Height = 351
Write-Host "Top:" $co.Top
Write-Host "Left:" $co.Left
Write-Host "Width:" $co.Width
Write-Host "Height:" $co.Height
